$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 03:52"

# Row 74: Nueva Zelanda - refreshed counts
$ws.Range("A74").Value = "Nueva Zelanda"
$ws.Range("B74").Value = 1461
$ws.Range("C74").Value = 5
$ws.Range("D74").Value = 1118
$ws.Range("E74").Value = 325
$ws.Range("F74").Value = 1
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 18

# Row 91: Bolivia moves up (new data)
$ws.Range("A91").Value = "Bolivia"
$ws.Range("B91").Value = 807
$ws.Range("C91").Value = 104
$ws.Range("D91").Value = 54
$ws.Range("E91").Value = 709
$ws.Range("F91").Value = 3
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 44

# Row 92: Republica de Chipre shifts down one slot
$ws.Range("A92").Value = "Republica de Chipre"
$ws.Range("B92").Value = 804
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 98
$ws.Range("E92").Value = 692
$ws.Range("F92").Value = 15
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 14

# Row 93: Letonia shifts down one slot
$ws.Range("A93").Value = "Letonia"
$ws.Range("B93").Value = 784
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 267
$ws.Range("E93").Value = 505
$ws.Range("F93").Value = 6
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 12

# Row 94: Principado de Andorra shifts down one slot
$ws.Range("A94").Value = "Principado de Andorra"
$ws.Range("B94").Value = 731
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 344
$ws.Range("E94").Value = 347
$ws.Range("F94").Value = 17
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 40

# Row 95: Crucero shifts down one slot
$ws.Range("A95").Value = "Crucero"
$ws.Range("B95").Value = 712
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 645
$ws.Range("E95").Value = 54
$ws.Range("F95").Value = 4
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 13

# Row 110: Guatemala moves up (new data)
$ws.Range("A110").Value = "Guatemala"
$ws.Range("B110").Value = 430
$ws.Range("C110").Value = 46
$ws.Range("D110").Value = 30
$ws.Range("E110").Value = 389
$ws.Range("F110").Value = 5
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 11

# Row 111: Taiwan shifts down one slot
$ws.Range("A111").Value = "Taiwan"
$ws.Range("B111").Value = 428
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 264
$ws.Range("E111").Value = 158
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 6

# Row 112: Sri Lanka shifts down one slot
$ws.Range("A112").Value = "Sri Lanka"
$ws.Range("B112").Value = 417
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 109
$ws.Range("E112").Value = 301
$ws.Range("F112").Value = 2
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 7

# Row 113: Reunion shifts down one slot
$ws.Range("A113").Value = "Reunion"
$ws.Range("B113").Value = 412
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 300
$ws.Range("E113").Value = 112
$ws.Range("F113").Value = 2
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 0

# Row 114: Consejo Danes para los Refugiados shifts down one slot
$ws.Range("A114").Value = "Consejo Danes para los Refugiados"
$ws.Range("B114").Value = 394
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 48
$ws.Range("E114").Value = 321
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 25

# Row 188: Dominica swaps up with Republica de Africa Central
$ws.Range("A188").Value = "Dominica"

# Row 189: Republica de Africa Central swaps down with Dominica
$ws.Range("A189").Value = "Republica de Africa Central"

# Row 197: Burundi swaps up with Islas Turcas y Caicos
$ws.Range("A197").Value = "Burundi"

# Row 198: Islas Turcas y Caicos swaps down with Burundi
$ws.Range("A198").Value = "Islas Turcas y Caicos"

# Row 210: Sudan del Sur swaps up with Bonaire, San Eustaquio y Saba
$ws.Range("A210").Value = "Sudan del Sur"

# Row 211: Bonaire, San Eustaquio y Saba swaps down with Sudan del Sur
$ws.Range("A211").Value = "Bonaire, San Eustaquio y Saba"
